$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.302
$ws.Range("E4").Value = 0.16
$ws.Range("G4").Value = 0.112
$ws.Range("H4").Value = 0.2
$ws.Range("J4").Value = 0.101
$ws.Range("K4").Value = 0.352
$ws.Range("L4").Value = 0.105
$ws.Range("M4").Value = 0.323
$ws.Range("N4").Value = 0.277
$ws.Range("O4").Value = 0.019
$ws.Range("P4").Value = 0.139
$ws.Range("Q4").Value = 0.542
$ws.Range("R4").Value = 0.214
$ws.Range("S4").Value = 0.462
$ws.Range("T4").Value = 0.299
$ws.Range("W4").Value = 0.239
$ws.Range("Y4").Value = 0.212
$ws.Range("Z4").Value = 0.467
$ws.Range("AA4").Value = 0.134
$ws.Range("AB4").Value = 0.365
$ws.Range("AE4").Value = 0.077
$ws.Range("AF4").Value = 0.722
$ws.Range("AG4").Value = 0.095
$ws.Range("AI4").Value = 0.639
$ws.Range("AJ4").Value = 0.175
$ws.Range("AK4").Value = 0.419
$ws.Range("AL4").Value = 0.6860000000000001
$ws.Range("AO4").Value = 0.6820000000000001
$ws.Range("B5").Value = 0.806
$ws.Range("C5").Value = 0.157
$ws.Range("D5").Value = 0.396
$ws.Range("E5").Value = 0.694
$ws.Range("F5").Value = 0.212
$ws.Range("G5").Value = 0.461
$ws.Range("H5").Value = 0.861
$ws.Range("I5").Value = 0.12
$ws.Range("J5").Value = 0.346
$ws.Range("K5").Value = 0.667
$ws.Range("L5").Value = 0.222
$ws.Range("M5").Value = 0.471
$ws.Range("N5").Value = 0.861
$ws.Range("O5").Value = 0.12
$ws.Range("P5").Value = 0.346
$ws.Range("Q5").Value = 0.611
$ws.Range("R5").Value = 0.238
$ws.Range("S5").Value = 0.487
$ws.Range("T5").Value = 0.611
$ws.Range("U5").Value = 0.238
$ws.Range("V5").Value = 0.487
$ws.Range("W5").Value = 0.722
$ws.Range("X5").Value = 0.201
$ws.Range("Y5").Value = 0.448
$ws.Range("Z5").Value = 0.833
$ws.Range("AA5").Value = 0.139
$ws.Range("AB5").Value = 0.373
$ws.Range("AC5").Value = 0.778
$ws.Range("AD5").Value = 0.173
$ws.Range("AE5").Value = 0.416
$ws.Range("AF5").Value = 0.972
$ws.Range("AG5").Value = 0.027
$ws.Range("AH5").Value = 0.164
$ws.Range("AI5").Value = 0.75
$ws.Range("AJ5").Value = 0.188
$ws.Range("AK5").Value = 0.433
$ws.Range("AL5").Value = 0.917
$ws.Range("AM5").Value = 0.076
$ws.Range("AN5").Value = 0.276
$ws.Range("AO5").Value = 0.88
$ws.Range("B6").Value = 0.439
$ws.Range("E6").Value = 0.26
$ws.Range("H6").Value = 0.325
$ws.Range("K6").Value = 0.461
$ws.Range("N6").Value = 0.419
$ws.Range("Q6").Value = 0.574
$ws.Range("T6").Value = 0.402
$ws.Range("W6").Value = 0.359
$ws.Range("Z6").Value = 0.598
$ws.Range("AF6").Value = 0.829
$ws.Range("AI6").Value = 0.6899999999999999
$ws.Range("AL6").Value = 0.785
$ws.Range("AO6").Value = 0.768
$ws.Range("B7").Value = 0.604
$ws.Range("E7").Value = 0.416
$ws.Range("H7").Value = 0.518
$ws.Range("K7").Value = 0.5659999999999999
$ws.Range("N7").Value = 0.606
$ws.Range("Q7").Value = 0.596
$ws.Range("T7").Value = 0.506
$ws.Range("W7").Value = 0.514
$ws.Range("Z7").Value = 0.72
$ws.Range("AC7").Value = 0.388
$ws.Range("AF7").Value = 0.909
$ws.Range("AI7").Value = 0.725
$ws.Range("AL7").Value = 0.859
$ws.Range("AO7").Value = 0.831
$ws.Range("B8").Value = 0.761
$ws.Range("C8").Value = 0.156
$ws.Range("D8").Value = 0.394
$ws.Range("E8").Value = 0.578
$ws.Range("H8").Value = 0.742
$ws.Range("I8").Value = 0.131
$ws.Range("J8").Value = 0.363
$ws.Range("K8").Value = 0.591
$ws.Range("L8").Value = 0.201
$ws.Range("M8").Value = 0.448
$ws.Range("N8").Value = 0.777
$ws.Range("O8").Value = 0.127
$ws.Range("P8").Value = 0.356
$ws.Range("Q8").Value = 0.58
$ws.Range("R8").Value = 0.224
$ws.Range("S8").Value = 0.473
$ws.Range("T8").Value = 0.528
$ws.Range("U8").Value = 0.202
$ws.Range("V8").Value = 0.45
$ws.Range("W8").Value = 0.653
$ws.Range("X8").Value = 0.187
$ws.Range("Y8").Value = 0.432
$ws.Range("Z8").Value = 0.765
$ws.Range("AA8").Value = 0.14
$ws.Range("AB8").Value = 0.375
$ws.Range("AC8").Value = 0.674
$ws.Range("AD8").Value = 0.174
$ws.Range("AE8").Value = 0.417
$ws.Range("AF8").Value = 0.887
$ws.Range("AG8").Value = 0.048
$ws.Range("AH8").Value = 0.22
$ws.Range("AI8").Value = 0.74
$ws.Range("AJ8").Value = 0.186
$ws.Range("AK8").Value = 0.431
$ws.Range("AL8").Value = 0.886
$ws.Range("AM8").Value = 0.082
$ws.Range("AN8").Value = 0.286
$ws.Range("AO8").Value = 0.838
$ws.Range("B9").Value = 0.694
$ws.Range("C9").Value = 0.212
$ws.Range("D9").Value = 0.461
$ws.Range("E9").Value = 0.444
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("H9").Value = 0.611
$ws.Range("I9").Value = 0.238
$ws.Range("J9").Value = 0.487
$ws.Range("K9").Value = 0.5
$ws.Range("N9").Value = 0.667
$ws.Range("O9").Value = 0.222
$ws.Range("P9").Value = 0.471
$ws.Range("Q9").Value = 0.528
$ws.Range("R9").Value = 0.249
$ws.Range("S9").Value = 0.499
$ws.Range("T9").Value = 0.417
$ws.Range("U9").Value = 0.243
$ws.Range("V9").Value = 0.493
$ws.Range("W9").Value = 0.556
$ws.Range("X9").Value = 0.247
$ws.Range("Y9").Value = 0.497
$ws.Range("Z9").Value = 0.667
$ws.Range("AA9").Value = 0.222
$ws.Range("AB9").Value = 0.471
$ws.Range("AC9").Value = 0.583
$ws.Range("AD9").Value = 0.243
$ws.Range("AE9").Value = 0.493
$ws.Range("AF9").Value = 0.75
$ws.Range("AG9").Value = 0.188
$ws.Range("AH9").Value = 0.433
$ws.Range("AI9").Value = 0.722
$ws.Range("AJ9").Value = 0.201
$ws.Range("AK9").Value = 0.448
$ws.Range("AL9").Value = 0.833
$ws.Range("AM9").Value = 0.139
$ws.Range("AN9").Value = 0.373
$ws.Range("AO9").Value = 0.768
$ws.Range("B10").Value = 0.806
$ws.Range("C10").Value = 0.157
$ws.Range("D10").Value = 0.396
$ws.Range("E10").Value = 0.611
$ws.Range("F10").Value = 0.238
$ws.Range("G10").Value = 0.487
$ws.Range("H10").Value = 0.778
$ws.Range("I10").Value = 0.173
$ws.Range("J10").Value = 0.416
$ws.Range("K10").Value = 0.667
$ws.Range("L10").Value = 0.222
$ws.Range("M10").Value = 0.471
$ws.Range("N10").Value = 0.833
$ws.Range("O10").Value = 0.139
$ws.Range("P10").Value = 0.373
$ws.Range("Q10").Value = 0.611
$ws.Range("R10").Value = 0.238
$ws.Range("S10").Value = 0.487
$ws.Range("T10").Value = 0.611
$ws.Range("U10").Value = 0.238
$ws.Range("V10").Value = 0.487
$ws.Range("W10").Value = 0.722
$ws.Range("X10").Value = 0.201
$ws.Range("Y10").Value = 0.448
$ws.Range("Z10").Value = 0.833
$ws.Range("AA10").Value = 0.139
$ws.Range("AB10").Value = 0.373
$ws.Range("AC10").Value = 0.667
$ws.Range("AD10").Value = 0.222
$ws.Range("AE10").Value = 0.471
$ws.Range("AF10").Value = 0.972
$ws.Range("AG10").Value = 0.027
$ws.Range("AH10").Value = 0.164
$ws.Range("AI10").Value = 0.75
$ws.Range("AJ10").Value = 0.188
$ws.Range("AK10").Value = 0.433
$ws.Range("AL10").Value = 0.917
$ws.Range("AM10").Value = 0.076
$ws.Range("AN10").Value = 0.276
$ws.Range("AO10").Value = 0.88
$ws.Range("B11").Value = 0.806
$ws.Range("C11").Value = 0.157
$ws.Range("D11").Value = 0.396
$ws.Range("E11").Value = 0.694
$ws.Range("F11").Value = 0.212
$ws.Range("G11").Value = 0.461
$ws.Range("H11").Value = 0.861
$ws.Range("I11").Value = 0.12
$ws.Range("J11").Value = 0.346
$ws.Range("K11").Value = 0.667
$ws.Range("L11").Value = 0.222
$ws.Range("M11").Value = 0.471
$ws.Range("N11").Value = 0.861
$ws.Range("O11").Value = 0.12
$ws.Range("P11").Value = 0.346
$ws.Range("Q11").Value = 0.611
$ws.Range("R11").Value = 0.238
$ws.Range("S11").Value = 0.487
$ws.Range("T11").Value = 0.611
$ws.Range("U11").Value = 0.238
$ws.Range("V11").Value = 0.487
$ws.Range("W11").Value = 0.722
$ws.Range("X11").Value = 0.201
$ws.Range("Y11").Value = 0.448
$ws.Range("Z11").Value = 0.833
$ws.Range("AA11").Value = 0.139
$ws.Range("AB11").Value = 0.373
$ws.Range("AC11").Value = 0.722
$ws.Range("AD11").Value = 0.201
$ws.Range("AE11").Value = 0.448
$ws.Range("AF11").Value = 0.972
$ws.Range("AG11").Value = 0.027
$ws.Range("AH11").Value = 0.164
$ws.Range("AI11").Value = 0.75
$ws.Range("AJ11").Value = 0.188
$ws.Range("AK11").Value = 0.433
$ws.Range("AL11").Value = 0.917
$ws.Range("AM11").Value = 0.076
$ws.Range("AN11").Value = 0.276
$ws.Range("AO11").Value = 0.88
$ws.Range("B12").Value = 1.172
$ws.Range("C12").Value = 0.212
$ws.Range("D12").Value = 0.46
$ws.Range("E12").Value = 1.68
$ws.Range("F12").Value = 1.098
$ws.Range("G12").Value = 1.048
$ws.Range("H12").Value = 1.613
$ws.Range("I12").Value = 1.334
$ws.Range("J12").Value = 1.155
$ws.Range("K12").Value = 1.417
$ws.Range("L12").Value = 0.576
$ws.Range("M12").Value = 0.759
$ws.Range("N12").Value = 1.355
$ws.Range("O12").Value = 0.552
$ws.Range("P12").Value = 0.743
$ws.Range("Z12").Value = 1.267
$ws.Range("AA12").Value = 0.329
$ws.Range("AB12").Value = 0.573
$ws.Range("AC12").Value = 1.786
$ws.Range("AD12").Value = 2.526
$ws.Range("AE12").Value = 1.589
$ws.Range("AF12").Value = 1.257
$ws.Range("AG12").Value = 0.248
$ws.Range("AH12").Value = 0.498
$ws.Range("AI12").Value = 1.037
$ws.Range("AJ12").Value = 0.036
$ws.Range("AK12").Value = 0.189
$ws.Range("AL12").Value = 1.091
$ws.Range("AM12").Value = 0.083
$ws.Range("AN12").Value = 0.287
$ws.Range("AO12").Value = 1.128
$ws.Range("B13").Value = 3.389
$ws.Range("C13").Value = 1.404
$ws.Range("D13").Value = 1.185
$ws.Range("E13").Value = 4.594
$ws.Range("F13").Value = 0.429
$ws.Range("G13").Value = 0.655
$ws.Range("H13").Value = 4.611
$ws.Range("I13").Value = 0.627
$ws.Range("J13").Value = 0.792
$ws.Range("K13").Value = 2.281
$ws.Range("L13").Value = 0.577
$ws.Range("M13").Value = 0.76
$ws.Range("N13").Value = 3.25
$ws.Range("O13").Value = 0.743
$ws.Range("P13").Value = 0.862
$ws.Range("Z13").Value = 2.5
$ws.Range("AA13").Value = 2.956
$ws.Range("AB13").Value = 1.719
$ws.Range("AC13").Value = 6.314
$ws.Range("AD13").Value = 2.216
$ws.Range("AE13").Value = 1.488
$ws.Range("AF13").Value = 1.639
$ws.Range("AG13").Value = 0.731
$ws.Range("AH13").Value = 0.855
$ws.Range("AI13").Value = 1.306
$ws.Range("AJ13").Value = 0.379
$ws.Range("AK13").Value = 0.616
$ws.Range("AL13").Value = 1.611
$ws.Range("AM13").Value = 0.738
$ws.Range("AN13").Value = 0.859
$ws.Range("AO13").Value = 1.519
